# results_BSD_withT_10_10.xlsx -- "new results abou BSD without turns"
#
# Mirrors columns A:I of each data block into a second, new block in
# columns K:S (A -> K, B:I -> L:S) and fills in the freshly computed
# "without turns" numbers for the rows the author already had results
# for. Also tweaks a couple of cosmetic bits (blank-cell formatting,
# column width, active selection) to track the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Mirror the "A" (row-label) column into "K" for every block, for
#    both format (style) and value/type (string vs number) -- this is
#    what makes K1 a string cell pointing at the same shared string as
#    A1, K30 a plain number matching A30, etc.
# ---------------------------------------------------------------------
$labelRowRanges = @("A1:A9", "A15:A23", "A29:A34", "A40:A45")
foreach ($src in $labelRowRanges) {
    $dst = $src -replace "^A", "K"

    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Copy the B:I number formatting into L:S for every row that gets
#    new numbers written into it.
# ---------------------------------------------------------------------
$formatRowRanges = @("B1:I3", "B15:I17", "B29:I34", "B40:I45")
foreach ($src in $formatRowRanges) {
    $dst = $src -replace "^B", "L"
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) New "without turns" values for L:S, keyed by row number.
# ---------------------------------------------------------------------
$lsData = @{
  1  = @("5","10","15","20","25","30","35","40")
  2  = @("6.4000000000000001E-2","0.214","0.43","0.54600000000000004","0.67800000000000005","0.79800000000000004","0.88600000000000001","0.90600000000000003")
  3  = @("6.4000000000000001E-2","0.23400000000000001","0.44","0.65600000000000003","0.77400000000000002","0.92","0.93200000000000005","0.92800000000000005")
  15 = @("5","10","15","20","25","30","35","40")
  16 = @("7.0000000000000007E-2","0.24199999999999999","0.49","0.61199999999999999","0.72199999999999998","0.86799999999999999","0.93400000000000005","0.96799999999999997")
  17 = @("6.4000000000000001E-2","0.248","0.49199999999999999","0.70599999999999996","0.876","0.97","0.98199999999999998","0.98399999999999999")
  29 = @("5","10","15","20","25","30","35","40")
  30 = @("6.4000000000000001E-2","0.23400000000000001","0.44","0.65600000000000003","0.77400000000000002","0.92","0.93200000000000005","0.92800000000000005")
  31 = @("0.122","0.41199999999999998","0.68400000000000005","0.85","0.97199999999999998","1","1","1")
  32 = @("0.214","0.502","0.76200000000000001","0.90400000000000003","0.99399999999999999","1","1","1")
  33 = @("0.24399999999999999","0.56999999999999995","0.82399999999999995","0.93799999999999994","1","1","1","1")
  34 = @("0.27200000000000002","0.6","0.84799999999999998","0.95199999999999996","1","1","1","1")
  40 = @("5","10","15","20","25","30","35","40")
  41 = @("7.3999999999999996E-2","0.27","0.51600000000000001","0.72599999999999998","0.88800000000000001","0.97599999999999998","0.98799999999999999","0.98599999999999999")
  42 = @("9.6000000000000002E-2","0.31","0.54800000000000004","0.76800000000000002","0.89400000000000002","0.98599999999999999","0.996","0.996")
  43 = @("0.11799999999999999","0.37","0.61799999999999999","0.78","0.89400000000000002","0.98599999999999999","1","1")
  44 = @("0.18","0.41","0.64200000000000002","0.81599999999999995","0.90600000000000003","0.98799999999999999","1","1")
  45 = @("0.218","0.45800000000000002","0.65800000000000003","0.82199999999999995","0.92400000000000004","0.98799999999999999","1","1")
}

$cols = @("L","M","N","O","P","Q","R","S")
foreach ($row in $lsData.Keys) {
    $vals = $lsData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = [double]$vals[$i]
    }
}

# ---------------------------------------------------------------------
# 4) Small formatting/view touch-ups to match the saved workbook.
# ---------------------------------------------------------------------
# A10 loses its yellow/highlight fill (blank separator cell).
$ws.Range("A35").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# New column K is widened to fit the "distance, london" header.
$ws.Columns("K").ColumnWidth = 15.5

# Scroll/selection state left by the author after entering the data.
$ws.Range("L48").Select() | Out-Null
